$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing data to make room for a new first column ---
# Insert a new column before column A (old A/B -> new B/C), preserving widths.
# Rows 1-2 were already empty, so no row shift is required for the new header row.
$ws.Columns.Item(1).Insert()

# --- New column A width (Schematic Bezeichner), closest reachable to the authored 20.43 ---
$ws.Columns.Item(1).ColumnWidth = 19.6

# --- Header row + new data, entered in the same order the strings were authored in ---
$ws.Range("B2").Value = "Bezeichnung"
$ws.Range("A2").Value = "Schematic Bezeichner"
$ws.Range("C2").Value = "Mouser Link"
$ws.Range("A2:C2").Font.Bold = $true

$ws.Range("A5").Value = "Q5A, Q5B"
$ws.Range("C6").Value = "SJ1-3523N"
$ws.Range("C5").Value = "BSD223PH6327XTSA1"
$ws.Range("B6").Value = "Audio Buchse"
$ws.Range("C7").Value = "09HCP-470M-50"
$ws.Range("A7").Value = "L1"
$ws.Range("B7").Value = "Tiefpass Spule 47uH"
$ws.Range("A6").Value = "J1"
$ws.Range("B8").Value = "Tiefpass Kapazität"
$ws.Range("C8").Value = "smd-kondensator-2220-330nf-100v"
$ws.Range("A8").Value = "C5"
$ws.Range("C9").Value = "BSD235NH6327XT"
$ws.Range("B5").Value = "Inverter P MOS Pair nach Optokoppler "
$ws.Range("B9").Value = "Optokoppler Driver N MOS Pair"

# --- Hyperlinks on the "Mouser Link" column ---
$ws.Hyperlinks.Add($ws.Range("C3"), "https://de.farnell.com/on-semiconductor/mc78m12cdtrkg/linearer-spann-regler-12v-0-5a/dp/2822588RL")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://de.farnell.com/camdenboss/ctb0305-3/anschlussblock-3polig/dp/3882640")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.mouser.de/c/?q=SJ1-3523N")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://www.mouser.de/c/?q=09HCP-470M-50")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.mouser.de/c/?q=BSD223PH6327XTSA1")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.mouser.de/c/?q=smd-kondensator-2220-330nf-100v")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://www.mouser.de/c/?q=BSD235NH6327XT")

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection matches authored state ---
$ws.Range("C28").Select()
